# "remove column from alcohol data"
# The measurement sheet had a spurious extra column M; its data was a
# duplicate/incorrect series and column N held the real values. Deleting
# column M shifts N left into M, dropping the file from A1:N119 to A1:M119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column M entirely (shifts N -> M, etc.)
$ws.Range("M:M").Delete() | Out-Null

# Move the selection onto the (now last) column, row 1
$ws.Range("M1").Select() | Out-Null
